$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''51.776.27'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '''  +4.96%  '
$ws.Range('E2').Style = "Normal"
$ws.Range('D3').Value = '''2.761.88'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '''  +5.30%  '
$ws.Range('E3').Style = "Normal"
$ws.Range('E4').Value = '''  +0.06%  '
$ws.Range('E4').Style = "Normal"
$ws.Range('D5').Value = '''116.49'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '''  +3.94%  '
$ws.Range('E5').Style = "Normal"
$ws.Range('D6').Value = '''333.19'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '''  +3.02%  '
$ws.Range('E6').Style = "Normal"
$ws.Range('E8').Value = '''  +0.03%  '
$ws.Range('E8').Style = "Normal"
$ws.Range('E9').Value = '''  +6.29%  '
$ws.Range('E9').Style = "Normal"
$ws.Range('D10').Value = '''41.83'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '''  +5.06%  '
$ws.Range('E10').Style = "Normal"
$ws.Range('D11').Value = '''0.0860'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '''  +5.96%  '
$ws.Range('E11').Style = "Normal"
$ws.Range('D12').Value = '''20.20'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '''  +2.01%  '
$ws.Range('E12').Style = "Normal"
$ws.Range('E13').Value = '''  +1.99%  '
$ws.Range('E13').Style = "Normal"
$ws.Range('E14').Value = '''  +5.30%  '
$ws.Range('E14').Style = "Normal"
$ws.Range('D15').Value = '''3.194.85'
$ws.Range('D15').Style = "Normal"
$ws.Range('D16').Value = '''2.779.84'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '''  +5.71%  '
$ws.Range('E16').Style = "Normal"
$ws.Range('D17').Value = '''0.890'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '''  +3.58%  '
$ws.Range('E17').Style = "Normal"
$ws.Range('D18').Value = '''51.703.39'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '''  +5.02%  '
$ws.Range('E18').Style = "Normal"
$ws.Range('D19').Value = '''3.21'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '''  +5.46%  '
$ws.Range('E19').Style = "Normal"
$ws.Range('D20').Value = '''13.52'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '''  +4.33%  '
$ws.Range('E20').Style = "Normal"
$ws.Range('D21').Value = '''6.88'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '''  +2.66%  '
$ws.Range('E21').Style = "Normal"
$ws.Range('E22').Value = '''  +3.11%  '
$ws.Range('E22').Style = "Normal"
$ws.Range('D23').Value = '''278.49'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '''  +3.29%  '
$ws.Range('E23').Style = "Normal"
$ws.Range('D24').Value = '''69.64'
$ws.Range('D24').Style = "Normal"
$ws.Range('D25').Value = '''2.68'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '''  +5.98%  '
$ws.Range('E25').Style = "Normal"
$ws.Range('E26').Value = '''  +2.29%  '
$ws.Range('E26').Style = "Normal"
$ws.Range('E27').Value = '''  +0.10%  '
$ws.Range('E27').Style = "Normal"
$ws.Range('D28').Value = '''10.18'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '''  -1.35%  '
$ws.Range('E28').Style = "Normal"
$ws.Range('E30').Value = '''  +2.12%  '
$ws.Range('E30').Style = "Normal"
$ws.Range('D31').Value = '''35.06'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '''  +0.31%  '
$ws.Range('E31').Style = "Normal"
$ws.Range('D32').Value = '''50.03'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '''  +0.90%  '
$ws.Range('E32').Style = "Normal"
$ws.Range('E33').Value = '''  +1.71%  '
$ws.Range('E33').Style = "Normal"
$ws.Range('D34').Value = '''0.0824'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '''  +1.20%  '
$ws.Range('E34').Style = "Normal"
$ws.Range('E35').Value = '''  -0.02%  '
$ws.Range('E35').Style = "Normal"
$ws.Range('D36').Value = '''18.97'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '''  -0.17%  '
$ws.Range('E36').Style = "Normal"
$ws.Range('D37').Value = '''5.01'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '''  +2.14%  '
$ws.Range('E37').Style = "Normal"
$ws.Range('E38').Value = '''  +2.27%  '
$ws.Range('E38').Style = "Normal"
$ws.Range('D39').Value = '''3.23'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '''  +3.27%  '
$ws.Range('E39').Style = "Normal"
$ws.Range('D40').Value = '''0.0351'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '''  +9.67%  '
$ws.Range('E40').Style = "Normal"
$ws.Range('D41').Value = '''127.04'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '''  -0.27%  '
$ws.Range('E41').Style = "Normal"
$ws.Range('D42').Value = '''23.24'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '''  +3.97%  '
$ws.Range('E42').Style = "Normal"
$ws.Range('E43').Value = '''  +3.32%  '
$ws.Range('E43').Style = "Normal"
$ws.Range('E44').Value = '''  +7.62%  '
$ws.Range('E44').Style = "Normal"
$ws.Range('E45').Value = '''  +13.02%  '
$ws.Range('E45').Style = "Normal"
$ws.Range('D46').Value = '''2.089.21'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '''  +1.56%  '
$ws.Range('E46').Style = "Normal"
$ws.Range('E47').Value = '''  +3.21%  '
$ws.Range('E47').Style = "Normal"
$ws.Range('E48').Value = '''  +4.85%  '
$ws.Range('E48').Style = "Normal"
$ws.Range('E49').Value = '''  +6.52%  '
$ws.Range('E49').Style = "Normal"
$ws.Range('E50').Value = '''  +1.29%  '
$ws.Range('E50').Style = "Normal"
$ws.Range('D51').Value = '''59.95'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '''  +1.68%  '
$ws.Range('E51').Style = "Normal"
